$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Source data (rows 2-6) that gets duplicated twice below the existing rows
$names = @("Bai 1", "Bai 2", "Bai 3", "Bai 4", "Bai 5")
$links = @(
    "https://www.youtube.com/embed/TU87tYb6qVg",
    "https://www.youtube.com/embed/Snn0-Im3kUc",
    "https://www.youtube.com/embed/TU87tYb6qVg",
    "https://www.youtube.com/embed/TU87tYb6qVg",
    "https://www.youtube.com/embed/TU87tYb6qVg"
)

# Append two copies of rows 2-6 (IDs 1-5, names Bai 1..5, link hyperlinks)
# into rows 7-11 and 12-16
for ($block = 0; $block -lt 2; $block++) {
    for ($i = 0; $i -lt 5; $i++) {
        $row = 7 + ($block * 5) + $i
        $id = $i + 1

        $ws.Cells.Item($row, 1).Value = $id
        $ws.Cells.Item($row, 2).Value = $names[$i]

        $cell = $ws.Cells.Item($row, 3)
        $cell.Value = $links[$i]
        $ws.Hyperlinks.Add($cell, $links[$i]) | Out-Null
        $cell.Style = "Hyperlink"
    }
}

# Update selection to match the final state of the workbook
$ws.Range("A12:D16").Select()
